$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.095.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5219"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2603"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06313"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.627.15"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.414"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5566"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0$([char]0x2085)8191"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.089.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.727"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.21"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.172"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.418"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.395"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05878"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.258"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.438"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.651"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9840"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.761"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.392"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5660"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8555"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.713"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.028.37"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.788.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0$([char]0x2088)109"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.084"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05150"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4217"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.70%  "
